$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 7285
    4  = 5497
    5  = 79
    10 = 86
    11 = 105
    13 = 42
    15 = 282
    17 = 10
    19 = 40
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
